$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("duration") values were stored as the text "5 Days" for every
# course. Replace them with the numeric day-count per row.
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 9

# Update the active selection to match the saved view.
$ws.Range("C12").Select()
